# Generate Report for Handoff
#
# The localization-status report gained a new handoff/handback cycle for the
# "ea62baf9-d577-4308-af66-a669f5670645" file (row 6 on every sheet). This
# update refreshes the recorded date/time stamps for that row on the
# Overview sheet and on each of the per-language detail sheets.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest Handoff Date" column (D) for the ea62baf9... row.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D6").Value = "2016-25-11 12:25:29"

# zh-cn detail sheet: "Latest Handoff Datetime" column (E) for the same row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E6").Value = "2016-03-11 12:25:27"

# de-de detail sheet: "Latest Handoff Datetime" column (E) for the same row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E6").Value = "2016-03-11 12:25:29"
